$wb = $excel.ActiveWorkbook

# This workbook tracks a fund's quarterly holdings, one sheet per quarter
# plus a running "总计" (totals) summary sheet. We need to:
#   1. Add a brand-new sheet named "2022-Q1" (placed right before "总计")
#      containing the per-fund holdings for that quarter.
#   2. Insert a new top row in "总计" for "2022-Q1" (count=2, value=0.09),
#      pushing the existing quarters down and renumbering the index column.
#
# The existing "总计" sheet is repurposed to become "2022-Q1" (its header
# row / box formatting already matches the other quarter sheets), and a
# fresh sheet is appended right after it to hold the (still-named) "总计"
# summary table.

# --- Step 1: turn the current "总计" sheet into "2022-Q1" ---
$newSheet = $wb.Worksheets.Item("总计")
$newSheet.Name = "2022-Q1"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund codes (B) and the numeric-looking stat columns (D:G) are stored as
# plain text in this workbook (e.g. "006440" keeps its leading zero), so
# force a text format before writing the values - otherwise Excel
# auto-coerces them to numbers and the leading zero is lost. Clear the
# format again afterwards so the cells don't carry a stray custom style.
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "006440"
$newSheet.Range("C2").Value = "中信建投中证500指数增强A"
$newSheet.Range("D2").Value = "5.78"
$newSheet.Range("E2").Value = "94.71"
$newSheet.Range("F2").Value = "1.01"
$newSheet.Range("G2").Value = "0.0584"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006441"
$newSheet.Range("C3").Value = "中信建投中证500指数增强C"
$newSheet.Range("D3").Value = "3.11"
$newSheet.Range("E3").Value = "94.71"
$newSheet.Range("F3").Value = "1.01"
$newSheet.Range("G3").Value = "0.0314"
$newSheet.Range("H3").Value = 5

$newSheet.Range("B2:B3").ClearFormats()
$newSheet.Range("D2:G3").ClearFormats()

# --- Step 2: append a fresh "总计" sheet after "2022-Q1" with the totals ---
$totalSheet = $wb.Worksheets.Add($null, $newSheet)
$totalSheet.Name = "总计"

# Reuse "2022-Q1" sheet's header formatting (bold, centered, thin border)
# for the new header row, then overwrite the captions for this sheet.
$newSheet.Range("B1:D1").Copy($totalSheet.Range("B1:D1"))
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# Copy the index-column (A) styling from the same sheet's header row box
# so the "0,1,2,..." column reads consistently with its row-1 siblings.
$newSheet.Range("A2").Copy($totalSheet.Range("A2"))
$totalSheet.Range("A2:A4").ClearFormats()
$newSheet.Range("A2").Copy()
$totalSheet.Range("A2:A4").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.09

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.29

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 8
$totalSheet.Range("D4").Value = 2.55

# Match the page margins used by its sibling data sheets (points:
# 0.75in=54, 1in=72, 0.5in=36) rather than Excel's new-sheet default.
$totalSheet.PageSetup.LeftMargin = 54
$totalSheet.PageSetup.RightMargin = 54
$totalSheet.PageSetup.TopMargin = 72
$totalSheet.PageSetup.BottomMargin = 72
$totalSheet.PageSetup.HeaderMargin = 36
$totalSheet.PageSetup.FooterMargin = 36
